$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column I (rows 4-10) into the new column J,
# then overwrite the values for the 2020 column.
$ws.Range("I4:I10").Copy()
$ws.Range("J4:J10").PasteSpecial(-4122)

$ws.Range("J4").Value = 2020
$ws.Range("J5").Value = 370
$ws.Range("J6").Value = 5
$ws.Range("J7").Value = 5
$ws.Range("J8").Value = 20
$ws.Range("J9").Value = 19
$ws.Range("J10").Value = 73

# J3 only needs the bottom border used by the rest of row 3 (no font/
# alignment carried over, unlike I3's style).
$ws.Range("J3").Borders.Item(9).LineStyle = 1
$ws.Range("J3").Borders.Item(9).Weight = -4138

# I8 changes from the "-" placeholder text to an actual number.
$ws.Range("I8").Value = 42

# I9 and I10 values were revised upward.
$ws.Range("I9").Value = 30
$ws.Range("I10").Value = 62
